$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new 2022 header column (F1) as text, reusing the exact header style (bold+centered)
# from the adjoining header cell (E1), then write "2022" as a text formula and collapse it
# back down to a static value so no formula / extra style survives.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F1").Formula = "=""2022"""
$ws.Range("F1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$excel.CutCopyMode = $false

# Update existing values (B:E) with refreshed figures and add new F column values

# Row 2 - DK
$ws.Range("B2").Value = 13.62512651192116
$ws.Range("C2").Value = 14.01869158878505
$ws.Range("D2").Value = 14.41015089163237
$ws.Range("E2").Value = 13.58607951783805
$ws.Range("F2").Value = 13.75910054372869

# Row 3 - ES
$ws.Range("B3").Value = 5.966137459807074
$ws.Range("C3").Value = 7.184048480109471
$ws.Range("D3").Value = 2.547654436882318
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0

# Row 4 - FI
$ws.Range("B4").Value = 48.21894005212858
$ws.Range("C4").Value = 47.99981415230219
$ws.Range("D4").Value = 47.28610407328006
$ws.Range("E4").Value = 45.65205927332173
$ws.Range("F4").Value = 45.37589150005846

# Row 5 - LU
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 49.76771196283391
$ws.Range("F5").Value = 49.83548766157461

# Row 6 - LV
$ws.Range("B6").Value = 43.42524594706942
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0

# Row 7 - NL
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 54.95245452898037
$ws.Range("F7").Value = 60.73966448725332

# Row 8 - SE
$ws.Range("B8").Value = 2.527362816151457
$ws.Range("C8").Value = 3.379581225804629
$ws.Range("D8").Value = 3.621518485479235
$ws.Range("E8").Value = 4.730001771536243
$ws.Range("F8").Value = 5.931610986514432
